$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (row 6) with the next scheduled Metaculus median snapshot,
# mirroring the formatting (date-formatted numbers) of the preceding rows.
$ws.Range("A6").Value = 46012
$ws.Range("B6").Value = 48848

# Match the style (date number format) used by the existing data rows.
$ws.Range("A6:B6").NumberFormat = $ws.Range("A5:B5").NumberFormat
